$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Generated on" date ---
$ws.Range("B2").Value = "21/1/2019"

# --- Protect the data area from Excel's automatic number/date detection ---
# (values such as "100000.00" or "2019-01-25" would otherwise be parsed as
# numbers/dates instead of being stored as literal text, which is what the
# target workbook requires)
$dataRange = $ws.Range("B6:E12")
$dataRange.NumberFormat = "@"

# --- Row 6 ---
$ws.Range("B6").Value = "lpo/Dream uniforms/87790"
$ws.Range("C6").Value = "Dream uniforms"
$ws.Range("D6").Value = "2019-01-25"
$ws.Range("E6").Value = "100000.00"

# --- Row 7 ---
$ws.Range("B7").Value = "lpo/Dream uniforms/70615"
$ws.Range("C7").Value = "Dream uniforms"
$ws.Range("D7").Value = "2019-01-25"
$ws.Range("E7").Value = "55000.00"

# --- Row 8 ---
$ws.Range("B8").Value = "lpo/Dream uniforms/35440"
$ws.Range("C8").Value = "Dream uniforms"
$ws.Range("D8").Value = "2019-01-26"
$ws.Range("E8").Value = "40000.00"

# --- Row 9 ---
$ws.Range("B9").Value = "lpo/Dream uniforms/27595"
$ws.Range("C9").Value = "Dream uniforms"
$ws.Range("D9").Value = "2019-01-19"
$ws.Range("E9").Value = "30000.00"

# --- Row 10 ---
$ws.Range("B10").Value = "lpo/NeuralStack/67991"
$ws.Range("C10").Value = "NeuralStack"
$ws.Range("D10").Value = "2019-01-25"
$ws.Range("E10").Value = "100000.00"

# --- Row 11 ---
$ws.Range("B11").Value = "lpo/Dream uniforms/DU/PO/2019-001"
$ws.Range("C11").Value = "Dream uniforms"
$ws.Range("D11").Value = "2019-01-20"
$ws.Range("E11").Value = "5850000.00"

# --- Row 12: clear S.NO/C.P.O#/CUSTOMER, turn into a "Total Amount" row ---
$ws.Range("D12").Value = "Total Amount"
$ws.Range("E12").Value = "6175000.00"

# --- Restore the default (General) format for the data area ---
$dataRange.Style = "Normal"

# --- Fully remove the now-unused S.NO/C.P.O#/CUSTOMER cells on row 12 ---
$ws.Range("A12:C12").Clear()
